$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.846.86"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "'1.919.39"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'241.70"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "'0.2973"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "'0.06767"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "'1.901.63"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "'17.08"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "'0.07319"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "'5.183"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "'89.40"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "'0.6723"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "'30.814.82"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'0.000007983"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'13.55"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'2.162.58"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'5.248"
$ws.Range("E22").Value = "  +8.59%  "
$ws.Range("D23").Value = "'201.80"
$ws.Range("E23").Value = "  +7.50%  "
$ws.Range("D24").Value = "'6.275"
$ws.Range("E24").Value = "  +3.61%  "
$ws.Range("D25").Value = "'9.666"
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("D26").Value = "'160.60"
$ws.Range("D27").Value = "'18.90"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").Value = "'1.974"
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").Value = "'1.429"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").Value = "'4.357"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "'0.09202"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").Value = "'4.071"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "'0.05191"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'0.7488"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "'1.123"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "'2.722"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "'0.01859"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "'0.9266"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").Value = "'2.081"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "'0.4501"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").Value = "'72.99"
$ws.Range("E42").Value = "  +25.99%  "
$ws.Range("D43").Value = "'107.87"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").Value = "'5.925"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("D45").Value = "'1.010"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "'0.1395"
$ws.Range("E46").Value = "  +3.88%  "
$ws.Range("D47").Value = "'7.705"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "'36.49"
$ws.Range("E48").Value = "  +8.82%  "
$ws.Range("D49").Value = "'9.024"
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("D50").Value = "'0.05945"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("D51").Value = "'0.4066"
$ws.Range("E51").Value = "  +3.18%  "
